# Updates the "Price" (column D) and a couple of "Volume(1h)" (column E)
# values in the crypto symbol list, matching the GitHub Actions refresh of
# Fri Dec 23 12:46:36 UTC 2022.
#
# The numeric-looking prices are stored as *text* in the sheet, so each one
# is written with a leading apostrophe (Excel's "treat as text" escape) to
# keep Excel from re-interpreting the string as a Number cell and dropping
# the significant trailing digits (e.g. "245.80" -> 245.8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($sheet, [string]$addr, [string]$value)
    $sheet.Range($addr).Value = "'" + $value
}

Set-TextValue $ws "D2"  "245.80"
Set-TextValue $ws "D3"  "22.00"
Set-TextValue $ws "D4"  "5.415"
Set-TextValue $ws "D5"  "0.05849"
Set-TextValue $ws "D6"  "3.389"
Set-TextValue $ws "D7"  "6.353"
Set-TextValue $ws "D8"  "0.8131"
Set-TextValue $ws "D9"  "0.9891"
$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"
Set-TextValue $ws "D10" "0.1429"
Set-TextValue $ws "D11" "0.07535"
Set-TextValue $ws "D12" "0.03344"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws "D13" "0.03007"
Set-TextValue $ws "D14" "4.191"
Set-TextValue $ws "D15" "0.09404"
Set-TextValue $ws "D16" "0.001587"
Set-TextValue $ws "D17" "0.04821"
Set-TextValue $ws "D18" "0.0005889"
Set-TextValue $ws "D19" "0.006152"
Set-TextValue $ws "D20" "0.004114"
Set-TextValue $ws "D21" "0.0009995"
Set-TextValue $ws "D23" "3.705"
Set-TextValue $ws "D24" "2.224"
Set-TextValue $ws "D25" "0.3258"
Set-TextValue $ws "D27" "0.0001290"
Set-TextValue $ws "D40" "0.03866"
Set-TextValue $ws "D41" "0.1079"
Set-TextValue $ws "D42" "0.002410"
Set-TextValue $ws "D43" "0.003030"
Set-TextValue $ws "D45" "0.00005602"
Set-TextValue $ws "D47" "0.3714"
Set-TextValue $ws "D48" "0.1444"
